$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly graded homework/quiz scores added for a batch of students.
# Column P = "第10章作业" (Ch.10 homework), R = "第5章作业" (Ch.5 homework)

$ws.Range("P4").Value = 5
$ws.Range("R4").Value = 5

$ws.Range("P12").Value = 5
$ws.Range("R12").Value = 5

$ws.Range("P15").Value = 5
$ws.Range("R15").Value = 5

$ws.Range("R17").Value = 5

$ws.Range("P25").Value = 5
$ws.Range("R25").Value = 5

$ws.Range("P27").Value = 5

$ws.Range("R34").Value = 5

$ws.Range("P37").Value = 5
$ws.Range("R37").Value = 5

$ws.Range("P46").Value = 5
$ws.Range("R46").Value = 5

$ws.Range("P48").Value = 5
$ws.Range("R48").Value = 5

$ws.Range("P49").Value = 4.5
$ws.Range("R49").Value = 5

$ws.Range("P54").Value = 5
$ws.Range("R54").Value = 5

$ws.Range("R69").Value = 5

$ws.Range("R71").Value = 5

$ws.Range("P73").Value = 5
$ws.Range("R73").Value = 5

$ws.Range("P74").Value = 5
$ws.Range("R74").Value = 5

$ws.Range("P77").Value = 5
$ws.Range("R77").Value = 5

$ws.Range("R79").Value = 5

$ws.Range("P81").Value = 5
$ws.Range("R81").Value = 5

$ws.Range("P87").Value = 5
$ws.Range("R87").Value = 5

$ws.Range("P90").Value = 5
$ws.Range("R90").Value = 5

# Scroll the view down and leave the selection on R54, matching where the
# grader was working when the workbook was last saved.
$excel.ActiveWindow.ScrollRow = 39
$ws.Range("R54").Select()
